$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.77"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").Value = "'22.55"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").Value = "'5.382"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").Value = "'0.05694"
$ws.Range("D5").ClearFormats()

$ws.Range("D6").Value = "'3.400"
$ws.Range("D6").ClearFormats()

$ws.Range("D8").Value = "'0.8119"
$ws.Range("D8").ClearFormats()

$ws.Range("D9").Value = "'0.9351"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").Value = "'0.1421"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").Value = "'0.07435"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").Value = "'0.03067"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").Value = "'0.03015"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").Value = "'0.09370"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").Value = "'3.719"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "'0.001588"
$ws.Range("D16").ClearFormats()

$ws.Range("D17").Value = "'0.04755"
$ws.Range("D17").ClearFormats()

$ws.Range("D18").Value = "'0.01827"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "'0.0005788"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "18OneONEWorstin24h"

$ws.Range("D20").Value = "'0.006443"
$ws.Range("D20").ClearFormats()

$ws.Range("D21").Value = "'0.005003"
$ws.Range("D21").ClearFormats()

$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D23").ClearFormats()

$ws.Range("D24").Value = "'3.695"
$ws.Range("D24").ClearFormats()

$ws.Range("D25").Value = "'2.153"
$ws.Range("D25").ClearFormats()

$ws.Range("D26").Value = "'0.3252"
$ws.Range("D26").ClearFormats()

$ws.Range("D40").Value = "'0.04007"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").Value = "'0.1069"
$ws.Range("D41").ClearFormats()

$ws.Range("D42").Value = "'0.002710"
$ws.Range("D42").ClearFormats()

$ws.Range("D43").Value = "'0.002966"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "'0.007485"
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = "'0.00005888"
$ws.Range("D45").ClearFormats()

$ws.Range("D47").Value = "'0.4998"
$ws.Range("D47").ClearFormats()

$ws.Range("D48").Value = "'0.2145"
$ws.Range("D48").ClearFormats()

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").ClearFormats()
